# Rename the (only) worksheet from the Arabic "prompt" placeholder title
# to the real data-sheet title.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "البيانات"

# Add the floating summary/insights text box that sits to the right of the
# table (around column H, row 6) and lists the five analysis bullet points.
$shp = $ws.Shapes.AddTextbox(1, 483.45, 79.35, 450, 146.59)
$shp.Name = "مربع نص 1"

$line1 = " إجمالي قيمة المبيعات لكل منطقة"
$line2 = "أكثر فئة منتج مبيعًا من حيث الكمية"
$line3 = "المنتج الأكثر تحقيقًا للإيرادات"
$line4 = "متوسط سعر الوحدة المباعة لكل فئة منتج"
$line5 = "نسبة مساهمة كل منطقة في إجمالي المبيعات"

$tr = $shp.TextFrame2.TextRange
$tr.Text = $line1 + [char]10 + $line2 + [char]10 + $line3 + [char]10 + $line4 + [char]10 + $line5

# Light grey fill (theme "Background 2" darkened 15%), matching the
# original shape style.
$shp.Fill.ForeColor.RGB = 14277081

Write-Host "Added text box and renamed sheet"
